# Applies the 06_Database.pptx edit:
#   1. Every cached "date" placeholder (slides, slide layouts, the slide
#      master and the notes master) is bumped from 2019/5/23 to 2019/5/28.
#   2. Slide 2's subtitle placeholder first bullet is re-worded from
#      "This section discusses Views of Django." to
#      "This section discusses Database of Django." (typed as a mid-run
#      replacement of "discusses Views " -> "discusses Database ", which
#      is how the run ends up split into three pieces).

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes) {
    $phs = $shapes.Placeholders
    for ($j = 1; $j -le $phs.Count; $j++) {
        $ph = $phs.Item($j)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            if ($ph.HasTextFrame -and $ph.TextFrame.HasText) {
                $tr = $ph.TextFrame.TextRange
                if ($tr.Text -eq "2019/5/23") {
                    $tr.Text = "2019/5/28"
                }
            }
        }
    }
}

# 1a. Slides.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    Update-DatePlaceholders $s.Shapes
}

# 1b. Slide layouts (hanging off the one slide master).
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DatePlaceholders $layout.Shapes
}

# 1c. Slide master itself.
Update-DatePlaceholders $p.SlideMaster.Shapes

# 1d. Notes master.
Update-DatePlaceholders $p.NotesMaster.Shapes

# 2. Slide 2 subtitle wording tweak.
$slide2 = $p.Slides.Item(2)
$subtitle = $slide2.Shapes.Item(2)
$firstPara = $subtitle.TextFrame.TextRange.Paragraphs(1, 1)
$target = $firstPara.Characters(14, 16)
if ($target.Text -eq "discusses Views ") {
    $target.Text = "discusses Database "
}
